# This document reorganizes the "body" text under each (fixed-position)
# heading of the LOB1003 syllabus:
#   - Objetivos' PT/EN blurb trade places with Programa resumido's PT/EN blurb.
#   - Docente(s) bullet <-> Avaliação's "NF=A avaliação..." sentence <->
#     Programa's PT body <-> Programa resumido's PT body <-> Bibliografia's
#     reference list <-> Docente(s) bullet (a 4-link rotation).
#   - Método/Critério/Norma de recuperação values shift up by one slot, and
#     the bibliography list slides in under "Norma de recuperação:".
# The heading/style sequence of the document never changes; only the runs
# inside each paragraph move around. We therefore read each paragraph's
# FormattedText (so runs, rPr, and <w:br/> breaks travel intact) before the
# paragraph that holds it gets overwritten, always processing a "read" before
# the matching "write" clobbers the source.

$d = $word.ActiveDocument

# --- Rotation: Docente bullet(9) -> Bibliografia body(19) ---------------
#     Objetivos body(6) -> Docente bullet(9)
#     Programa resumido body(11) -> Objetivos body(6)
#     Programa body(14) -> Programa resumido body(11)
$tmp19 = $d.Paragraphs.Item(9).Range.FormattedText
$d.Paragraphs.Item(19).Range.FormattedText = $tmp19

$tmp9 = $d.Paragraphs.Item(6).Range.FormattedText
$d.Paragraphs.Item(9).Range.FormattedText = $tmp9

$tmp6 = $d.Paragraphs.Item(11).Range.FormattedText
$d.Paragraphs.Item(6).Range.FormattedText = $tmp6

$tmp11 = $d.Paragraphs.Item(14).Range.FormattedText
$d.Paragraphs.Item(11).Range.FormattedText = $tmp11

# --- Swap: Objetivos EN blurb(7) <-> Programa resumido EN blurb(12) -----
$tmp7 = $d.Paragraphs.Item(7).Range.FormattedText
$tmp12 = $d.Paragraphs.Item(12).Range.FormattedText
$d.Paragraphs.Item(7).Range.FormattedText = $tmp12
$d.Paragraphs.Item(12).Range.FormattedText = $tmp7

# --- Programa's body(14) is now free (its old text was already captured
#     above as $tmp11); fill it with the "NF=A avaliação..." sentence that
#     used to live inside the Avaliação bullet list (paragraph 17). --------
$d.Paragraphs.Item(14).Range.Text = "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n."

# --- Avaliação bullet list (paragraph 17): the three "value" runs shift up
#     by one (Método's value <- old Critério's value, Critério's value <-
#     old Norma's value), and the bibliography list slides in as the new
#     "Norma de recuperação:" value. Process back-to-front so every search
#     string is still unique in the paragraph at the moment it is matched.
$nl = [char]11
$bibliografia = "STEWART, James. Cálculo São Paulo: Cengage Learning, 2009. v.1." + $nl + $nl + `
    "ANTON, Howard. Cálculo: um novo horizonte. Porto Alegre: Bookman, 2007." + $nl + $nl + `
    "THOMAS, George B. Cálculo São Paulo: Pearson Addison  Wesley, 2009. v.1," + $nl + $nl + `
    "GUIDORIZZI, Hamilton. Um curso de cálculo. Rio de Janeiro: Livros Técnicos e Científicos, 2001. v.1." + $nl + $nl + `
    "FLEMMING, Diva M.; GONÇALVES, Mirian B. Cálculo A. São Paulo: Pearson Prentice Hall, 2009."

$p17a = $d.Paragraphs.Item(17).Range.Duplicate
[void]$p17a.Find.Execute(
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $bibliografia, 2)

$p17b = $d.Paragraphs.Item(17).Range.Duplicate
[void]$p17b.Find.Execute(
    "NF≥ 5,0.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "(NF+RC)/2 ≥ 5,0, onde RC é uma prova de recuperação a ser aplicada.", 2)

$p17c = $d.Paragraphs.Item(17).Range.Duplicate
[void]$p17c.Find.Execute(
    "NF=A avaliação será composta por provas, listas, projetos, seminários e outras formas que farão a composição das notas, sendo estipulada a média final a somatória destas notas (N), com no mínimo duas avaliações, sendo: (N1+...+Nn)/n.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "NF≥ 5,0.", 2)
